$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(102, 1).Value = "Data Content Type"
$ws.Cells.Item(102, 2).Value = "Epidemiologic"
$ws.Cells.Item(102, 3).Value = "Relating to the study of the distribution and determinants of health-related states or events (including disease) in populations, and the application of this study to the control of diseases and other health problems."

$ws.Cells.Item(103, 1).Value = "Data Element"
$ws.Cells.Item(103, 2).Value = "Proband"
$ws.Cells.Item(103, 3).Value = "A proband is a person in a family to receive genetic counseling and/or testing for a suspected hereditary risk or diagnosed disease. A proband may or may not be affected with the disease in question. If the value is true, then the case subject may have been diagnosed with the disease under studied. If the value is false, then the case subject is a member of the family of a proband study participant. The proband indicator for the case carries over to a sample taken from a case subject. (Definition based on https://www.cancer.gov/publications/dictionaries/genetics-dictionary/def/proband.)"

$ws.Range("C102:C103").WrapText = $true

$ws.Rows.Item(102).RowHeight = 29
$ws.Rows.Item(103).RowHeight = 87

$excel.ActiveWindow.ScrollRow = 94
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("C106").Select()
